# Add a new worksheet ("Sheet2") between "Second Sheet" and "Empty Sheet",
# containing a small new block of test data, matching the upstream commit
# "Add more test data to the Excel file".

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("Second Sheet")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)

# New table: header row at row 5 (Col1 in C5, Col3 in E5 - D5 left blank),
# two rows of numeric data below it.
$newSheet.Range("C5").Value = "Col1"
$newSheet.Range("E5").Value = "Col3"

$newSheet.Range("C6").Value = 34
$newSheet.Range("D6").Value = 35
$newSheet.Range("E6").Value = 23

$newSheet.Range("C7").Value = 14
$newSheet.Range("D7").Value = 27
$newSheet.Range("E7").Value = 69

# Match the author's final selection/active cell on the new sheet.
$newSheet.Range("C8").Select() | Out-Null
